$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
